$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-6 from 45212 to 45221
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45221
}
